$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the old header value in B1 ("Jun_13") before the column shift
$oldB1 = $ws.Cells.Item(1, 2).Text

# Insert two new columns at C (pushes the existing column C -> E)
$ws.Columns("C:D").Insert()

# Match the width/formatting of the neighbouring (old) column on the two
# newly inserted columns, same as Excel does when inserting columns.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14

# Header row: new "Jun_17" / "Jun_15" columns, then the old B1/C1 values
# shifted right into D1/E1.
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"
$ws.Cells.Item(1, 4).Value = $oldB1

# Fill the two new columns (C, D) for every data row with "UN", matching
# column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
